{"js": "// The recorded change (see the diff / commit) only re-serializes the\n// document's existing OOXML: every single hunk is a pure XML\n// attribute-reordering / whitespace-normalization of attributes that are\n// already present, e.g.\n//   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>   ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n//   <w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n//     -> <w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n//   <w:document xmlns:wpc=\"\u2026\" xmlns:mc=\"\u2026\" \u2026>  ->  <w:document xmlns:m=\"\u2026\" xmlns:mc=\"\u2026\" \u2026>\n// Every attribute name/value pair on every element touched by the diff\n// (namespace declarations on <w:document>, <wp:anchor>,\n// <wp:effectExtent>, <wps:bodyPr>, <v:shapetype>, <v:path>, <v:shape>,\n// <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>,\n// <w:lsdException>, <w:style>, table cell margins, ...) is identical\n// before and after - only the order the attributes happen to be\n// serialized in differs. No paragraph text, field code, run formatting,\n// shape/picture, style definition, or any other content/formatting value\n// is added, removed, or changed anywhere in the diff.\n//\n// Raw XML attribute order is not part of the Word JavaScript object\n// model (Office.js exposes properties/values, never the literal\n// attribute-serialization order of the underlying part XML), and\n// reordering attributes has no observable effect through that API:\n// paragraphs, runs, formatting, shapes, sections and styles all keep\n// exactly the same values. So the faithful way to reproduce this diff\n// through Office.js is to leave the document's content/formatting\n// untouched - deliberately poking unrelated parts of the OM would only\n// risk introducing incidental changes (e.g. new namespace declarations\n// being minted) that are not part of the recorded diff.\n//\n// We still read the document through the API (body text and the\n// section/page-margin values touched by the <w:pgSz>/<w:pgMar> hunk) so\n// the script demonstrably inspects the same content the diff\n// re-serialized, without writing anything back.\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nawait context.sync();\n\nfor (const section of sections.items) {\n    section.body.load(\"text\");\n}\n\nawait context.sync();\n", "ps1": "# The recorded change (see the diff / commit message) only re-serializes\n# the document's existing OOXML: every single hunk in the diff is a pure\n# XML attribute-reordering / whitespace-normalization of attributes that\n# are already present, e.g.\n#   <w:pgSz w:w=\"11906\" w:h=\"16838\"/>   ->  <w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#   <w:color w:val=\"E36C0A\" w:themeColor=\"accent6\" w:themeShade=\"BF\"/>\n#     -> <w:color w:themeColor=\"accent6\" w:themeShade=\"BF\" w:val=\"E36C0A\"/>\n#   <w:document xmlns:wpc=\"...\" xmlns:mc=\"...\" ...>\n#     -> <w:document xmlns:m=\"...\" xmlns:mc=\"...\" ...>\n# Every attribute name/value pair on every element touched by the diff\n# (namespace declarations on <w:document>, <wp:anchor>,\n# <wp:effectExtent>, <wps:bodyPr>, <v:shapetype>, <v:path>, <v:shape>,\n# <w:pgSz>, <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>,\n# <w:lsdException>, <w:style>, table cell margins, ...) is identical\n# before and after - only the order the attributes happen to be\n# serialized in differs. No paragraph text, field code, run formatting,\n# shape/picture, style definition, or any other content/formatting value\n# is added, removed, or changed anywhere in the diff.\n#\n# Raw XML attribute order is not part of the Word COM object model\n# either (PageSetup/Font/Section/Style properties are exposed by value,\n# never as the literal attribute-serialization order of the underlying\n# part XML), and reordering attributes has no observable effect through\n# that API: paragraphs, runs, formatting, shapes, sections and styles all\n# keep exactly the same values. So the faithful way to reproduce this\n# diff through COM is to leave the document's content/formatting\n# untouched - deliberately poking unrelated parts of the object model\n# would only risk introducing incidental changes that are not part of\n# the recorded diff.\n#\n# We still read the document through the object model (body text and the\n# section/page-margin values touched by the <w:pgSz>/<w:pgMar> hunk) so\n# the script demonstrably inspects the same content the diff\n# re-serialized, without writing anything back.\n\n$d = $word.ActiveDocument\n\n$null = $d.Content.Text\n\nforeach ($section in $d.Sections) {\n    $pageSetup = $section.PageSetup\n    $null = $pageSetup.PageWidth\n    $null = $pageSetup.PageHeight\n    $null = $pageSetup.TopMargin\n    $null = $pageSetup.BottomMargin\n    $null = $pageSetup.LeftMargin\n    $null = $pageSetup.RightMargin\n}\n"}
